# Update the Spp1-Itga9 LR-pairs sheet with newly computed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.3399353333333333
$ws.Range("H2").Value = 1.019806
$ws.Range("I2").Value = 0.09929991924017606
$ws.Range("J2").Value = 0.09929991924017606
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8366046666666667
$ws.Range("N2").Value = 2.509814
$ws.Range("O2").Value = 0.08025679986157715
$ws.Range("P2").Value = 0.08025679986157715
$ws.Range("Q2").Value = 0.2843914862315555
$ws.Range("R2").Value = 2.559523376084
$ws.Range("S2").Value = 0.007969493744729585
$ws.Range("T2").Value = 0.007969493744729585
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.3399353333333333
$ws.Range("H3").Value = 1.019806
$ws.Range("I3").Value = 0.09929991924017606
$ws.Range("J3").Value = 0.09929991924017606
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 7.939250333333333
$ws.Range("N3").Value = 23.817751
$ws.Range("O3").Value = 0.7616247559221037
$ws.Range("P3").Value = 0.7616247559221038
$ws.Range("Q3").Value = 2.698831708478444
$ws.Range("R3").Value = 24.289485376306
$ws.Range("S3").Value = 0.07562927675438369
$ws.Range("T3").Value = 0.0756292767543837
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.3399353333333333
$ws.Range("H4").Value = 1.019806
$ws.Range("I4").Value = 0.09929991924017606
$ws.Range("J4").Value = 0.09929991924017606
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.648242
$ws.Range("N4").Value = 4.944726
$ws.Range("O4").Value = 0.1581184442163192
$ws.Range("P4").Value = 0.1581184442163192
$ws.Range("Q4").Value = 0.560295693684
$ws.Range("R4").Value = 5.042661243156
$ws.Range("S4").Value = 0.01570114874106278
$ws.Range("T4").Value = 0.01570114874106278
$ws.Range("I5").Value = 0.4094685684206303
$ws.Range("J5").Value = 0.4094685684206303
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8366046666666667
$ws.Range("N5").Value = 2.509814
$ws.Range("O5").Value = 0.08025679986157715
$ws.Range("P5").Value = 0.08025679986157715
$ws.Range("Q5").Value = 1.172703619794444
$ws.Range("R5").Value = 10.55433257815
$ws.Range("S5").Value = 0.03286263694534104
$ws.Range("T5").Value = 0.03286263694534104
$ws.Range("I6").Value = 0.4094685684206303
$ws.Range("J6").Value = 0.4094685684206303
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 7.939250333333333
$ws.Range("N6").Value = 23.817751
$ws.Range("O6").Value = 0.7616247559221037
$ws.Range("P6").Value = 0.7616247559221038
$ws.Range("Q6").Value = 11.12877799433055
$ws.Range("R6").Value = 100.159001948975
$ws.Range("S6").Value = 0.3118613984811358
$ws.Range("T6").Value = 0.3118613984811358
$ws.Range("I7").Value = 0.4094685684206303
$ws.Range("J7").Value = 0.4094685684206303
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.648242
$ws.Range("N7").Value = 4.944726
$ws.Range("O7").Value = 0.1581184442163192
$ws.Range("P7").Value = 0.1581184442163192
$ws.Range("Q7").Value = 2.310409488149999
$ws.Range("R7").Value = 20.79368539335
$ws.Range("S7").Value = 0.0647445329941535
$ws.Range("T7").Value = 0.06474453299415353
$ws.Range("G8").Value = 1.681642333333333
$ws.Range("H8").Value = 5.044927
$ws.Range("I8").Value = 0.4912315123391937
$ws.Range("J8").Value = 0.4912315123391937
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.8366046666666667
$ws.Range("N8").Value = 2.509814
$ws.Range("O8").Value = 0.08025679986157715
$ws.Range("P8").Value = 0.08025679986157715
$ws.Range("Q8").Value = 1.406869823730889
$ws.Range("R8").Value = 12.661828413578
$ws.Range("S8").Value = 0.03942466917150653
$ws.Range("T8").Value = 0.03942466917150653
$ws.Range("G9").Value = 1.681642333333333
$ws.Range("H9").Value = 5.044927
$ws.Range("I9").Value = 0.4912315123391937
$ws.Range("J9").Value = 0.4912315123391937
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 7.939250333333333
$ws.Range("N9").Value = 23.817751
$ws.Range("O9").Value = 0.7616247559221037
$ws.Range("P9").Value = 0.7616247559221038
$ws.Range("Q9").Value = 13.35097945546411
$ws.Range("R9").Value = 120.158815099177
$ws.Range("S9").Value = 0.3741340806865842
$ws.Range("T9").Value = 0.3741340806865843
$ws.Range("G10").Value = 1.681642333333333
$ws.Range("H10").Value = 5.044927
$ws.Range("I10").Value = 0.4912315123391937
$ws.Range("J10").Value = 0.4912315123391937
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.648242
$ws.Range("N10").Value = 4.944726
$ws.Range("O10").Value = 0.1581184442163192
$ws.Range("P10").Value = 0.1581184442163192
$ws.Range("Q10").Value = 2.771753522778
$ws.Range("R10").Value = 24.945781705002
$ws.Range("S10").Value = 0.0776727624811029
$ws.Range("T10").Value = 0.07767276248110291
